$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($CellRef, $Val)
    $rng = $ws.Range($CellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $Val
    $rng.ClearFormats()
}

$ws.Range("D2").Value = "28.000.73"
$ws.Range("E2").Value = "  +1.93%  "
$ws.Range("D3").Value = "1.907.23"
Set-TextValue "D4" "1.003"
$ws.Range("E4").Value = "  -0.67%  "
Set-TextValue "D5" "317.22"
$ws.Range("E5").Value = "  +1.84%  "
Set-TextValue "D6" "1.003"
$ws.Range("E6").Value = "  -0.68%  "
Set-TextValue "D7" "0.4838"
$ws.Range("E7").Value = "  +1.35%  "
Set-TextValue "D8" "0.3801"
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +0.84%  "
Set-TextValue "D10" "0.9330"
$ws.Range("E10").Value = "  +0.14%  "
Set-TextValue "D11" "20.77"
$ws.Range("E11").Value = "  +0.14%  "
Set-TextValue "D12" "0.07748"
$ws.Range("D13").Value = "1.936.90"
$ws.Range("E13").Value = "  +3.86%  "
Set-TextValue "D14" "5.483"
$ws.Range("E14").Value = "  +0.68%  "
Set-TextValue "D15" "6.644"
$ws.Range("E15").Value = "  +1.44%  "
Set-TextValue "D16" "91.84"
$ws.Range("E16").Value = "  +1.76%  "
$ws.Range("E17").Value = "  -0.69%  "
Set-TextValue "D18" "0.000008879"
$ws.Range("E18").Value = "  +0.69%  "
$ws.Range("E19").Value = "  -0.59%  "
$ws.Range("D20").Value = "28.033.50"
$ws.Range("E20").Value = "  +1.86%  "
$ws.Range("E21").Value = "  +0.32%  "
Set-TextValue "D22" "5.148"
$ws.Range("E22").Value = "  +1.03%  "
$ws.Range("D23").Value = "2.153.00"
$ws.Range("E23").Value = "  +1.82%  "
Set-TextValue "D24" "10.89"
$ws.Range("E24").Value = "  +1.96%  "
Set-TextValue "D25" "156.14"
$ws.Range("E25").Value = "  +0.74%  "
Set-TextValue "D26" "1.918"
$ws.Range("E26").Value = "  -1.27%  "
$ws.Range("E27").Value = "  +0.24%  "
Set-TextValue "D28" "2.121"
$ws.Range("E28").Value = "  +5.74%  "
Set-TextValue "D29" "117.28"
$ws.Range("E29").Value = "  +1.61%  "
Set-TextValue "D30" "4.980"
$ws.Range("E30").Value = "  +0.54%  "
Set-TextValue "D31" "0.08940"
$ws.Range("E31").Value = "  +0.49%  "
Set-TextValue "D32" "3.269"
$ws.Range("E32").Value = "  -1.91%  "
Set-TextValue "D33" "1.253"
$ws.Range("E33").Value = "  +4.05%  "
Set-TextValue "D34" "0.7700"
$ws.Range("E34").Value = "  +2.12%  "
Set-TextValue "D35" "4.667"
$ws.Range("E35").Value = "  +1.73%  "
Set-TextValue "D36" "2.584"
$ws.Range("E36").Value = "  -4.18%  "
Set-TextValue "D37" "0.02057"
$ws.Range("E37").Value = "  +0.86%  "
Set-TextValue "D39" "0.5499"
$ws.Range("E39").Value = "  -1.16%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D40" "3.002"
$ws.Range("E40").Value = "  +0.52%  "
$ws.Range("B41").Value = "Hedera"
$ws.Range("C41").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D41" "0.05274"
$ws.Range("E41").Value = "  -0.08%  "
Set-TextValue "D42" "6.931"
$ws.Range("E42").Value = "  -1.42%  "
$ws.Range("E43").Value = "  +0.77%  "
Set-TextValue "D44" "8.503"
$ws.Range("E44").Value = "  -0.98%  "
Set-TextValue "D45" "110.88"
$ws.Range("E45").Value = "  +7.65%  "
$ws.Range("E46").Value = "  +0.19%  "
Set-TextValue "D47" "0.4819"
$ws.Range("E47").Value = "  -1.02%  "
$ws.Range("E48").Value = "  -0.73%  "
Set-TextValue "D49" "1.645"
$ws.Range("E49").Value = "  -1.13%  "
Set-TextValue "D50" "68.01"
$ws.Range("E50").Value = "  +0.84%  "
Set-TextValue "D51" "0.06073"
$ws.Range("E51").Value = "  -0.20%  "
